$wb = $excel.ActiveWorkbook

# "IT" sheet (sheetId 2 / second tab) becomes the active sheet/tab.
$ws = $wb.Worksheets.Item("IT")
$ws.Activate()

# Update the year value in B2 from 2019 to 2020.
$ws.Range("B2").Value = 2020

# The active cell/selection on the IT sheet becomes B3.
$ws.Range("B3").Select()
